$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 411.42856
$ws.Range("I12").Value = 193.33333
$ws.Range("K12").Value = 193.33333
$ws.Range("M12").Value = -23.33332999999999
$ws.Range("H32").Value = 4371.4287
$ws.Range("I32").Value = 4028.5715
$ws.Range("K32").Value = 4028.5715
$ws.Range("M32").Value = -3702.5715
$ws.Range("H43").Value = 8139.08
$ws.Range("I43").Value = 8087.2354
$ws.Range("K43").Value = 8087.2354
$ws.Range("M43").Value = -8018.2354
$ws.Range("H53").Value = 455.53845
$ws.Range("I53").Value = 466.125
$ws.Range("J53").Value = 438.6
$ws.Range("K53").Value = 466.125
$ws.Range("L53").Value = 438.6
$ws.Range("M53").Value = 170.875
$ws.Range("N53").Value = -1712.6
$ws.Range("H98").Value = 58826228
$ws.Range("I98").Value = 71431010
$ws.Range("K98").Value = 71431010
$ws.Range("M98").Value = -71429512
$ws.Range("H106").Value = 13033.583
$ws.Range("I106").Value = 3601.25
$ws.Range("K106").Value = 3601.25
$ws.Range("M106").Value = -2970.25
$ws.Range("H116").Value = 6693
$ws.Range("I116").Value = 6737.875
$ws.Range("K116").Value = 6737.875
$ws.Range("M116").Value = -3295.875
$ws.Range("H122").Value = 58826228
$ws.Range("I122").Value = 71431010
$ws.Range("K122").Value = 214293030
$ws.Range("M122").Value = -214290580
$ws.Range("H124").Value = 106663.336
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 106663.336
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 106663.336
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -116483.336
$ws.Range("H132").Value = 1212.674
$ws.Range("I132").Value = 740.1795
$ws.Range("K132").Value = 2220.5385
$ws.Range("M132").Value = 309.4615000000003
$ws.Range("H133").Value = 79999.336
$ws.Range("J133").Value = 79999.336
$ws.Range("L133").Value = 79999.336
$ws.Range("N133").Value = -90119.336
$ws.Range("H137").Value = 8487.933999999999
$ws.Range("I137").Value = 8022.9287
$ws.Range("J137").Value = 14998
$ws.Range("K137").Value = 24068.7861
$ws.Range("L137").Value = 44994
$ws.Range("M137").Value = -21518.7861
$ws.Range("N137").Value = -50094
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 92714.5
$ws.Range("J135").Value = 92714.5
$ws.Range("L135").Value = 92714.5
$ws.Range("N135").Value = -102854.5
$ws.Range("H137").Value = 69996.5
$ws.Range("J137").Value = 69996.5
$ws.Range("L137").Value = 69996.5
$ws.Range("N137").Value = -80196.5
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1169215.2
$ws.Range("I31").Value = 3249.889
$ws.Range("K31").Value = 3249.889
$ws.Range("M31").Value = -2954.889
$ws.Range("H34").Value = 1169215.2
$ws.Range("I34").Value = 3249.889
$ws.Range("K34").Value = 3249.889
$ws.Range("M34").Value = -3047.889
$ws.Range("H99").Value = 2967.5
$ws.Range("I99").Value = 2811
$ws.Range("K99").Value = 2811
$ws.Range("M99").Value = -1313
$ws.Range("H124").Value = 42160.75
$ws.Range("J124").Value = 42160.75
$ws.Range("L124").Value = 42160.75
$ws.Range("N124").Value = -47070.75
$ws.Range("H126").Value = 2967.5
$ws.Range("I126").Value = 2811
$ws.Range("K126").Value = 8433
$ws.Range("M126").Value = -5963
$ws.Range("H127").Value = 67395.60000000001
$ws.Range("J127").Value = 67395.60000000001
$ws.Range("L127").Value = 67395.60000000001
$ws.Range("N127").Value = -77315.60000000001
$ws.Range("H132").Value = 6100.9473
$ws.Range("I132").Value = 3493.077
$ws.Range("J132").Value = 11751.333
$ws.Range("K132").Value = 10479.231
$ws.Range("L132").Value = 35253.999
$ws.Range("M132").Value = -7949.231
$ws.Range("N132").Value = -40313.999
$ws.Range("H134").Value = 3229.1052
$ws.Range("I134").Value = 2647.0625
$ws.Range("K134").Value = 7941.1875
$ws.Range("M134").Value = -5406.1875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 67249.5
$ws.Range("J37").Value = 67249.5
$ws.Range("L37").Value = 201748.5
$ws.Range("N37").Value = -201972.5
$ws.Range("H40").Value = 339.57144
$ws.Range("I40").Value = 339.57144
$ws.Range("K40").Value = 1358.28576
$ws.Range("M40").Value = -1289.28576
$ws.Range("H134").Value = 9560.25
$ws.Range("J134").Value = 12775.862
$ws.Range("L134").Value = 38327.586
$ws.Range("N134").Value = -48467.586
$ws.Range("H136").Value = 10377.2
$ws.Range("J136").Value = 14999.5
$ws.Range("L136").Value = 44998.5
$ws.Range("N136").Value = -55198.5
$ws.Range("H141").Value = 241402.23
$ws.Range("I141").Value = 505621.5
$ws.Range("J141").Value = 14928.571
$ws.Range("K141").Value = 1516864.5
$ws.Range("L141").Value = 44785.713
$ws.Range("M141").Value = -1511684.5
$ws.Range("N141").Value = -55145.713
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4658.222
$ws.Range("I126").Value = 4231
$ws.Range("K126").Value = 12693
$ws.Range("M126").Value = -10223
$ws.Range("H134").Value = 129999.5
$ws.Range("J134").Value = 129999.5
$ws.Range("L134").Value = 389998.5
$ws.Range("N134").Value = -395068.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10838.571
$ws.Range("I7").Value = 4357.4546
$ws.Range("K7").Value = 4357.4546
$ws.Range("M7").Value = -4245.4546
$ws.Range("H22").Value = 2590.2
$ws.Range("I22").Value = 2557.1428
$ws.Range("J22").Value = 2667.3333
$ws.Range("K22").Value = 2557.1428
$ws.Range("L22").Value = 2667.3333
$ws.Range("M22").Value = -2262.1428
$ws.Range("N22").Value = -3257.3333
$ws.Range("H27").Value = 2590.2
$ws.Range("I27").Value = 2557.1428
$ws.Range("J27").Value = 2667.3333
$ws.Range("K27").Value = 2557.1428
$ws.Range("L27").Value = 2667.3333
$ws.Range("M27").Value = -2450.1428
$ws.Range("N27").Value = -2881.3333
$ws.Range("H40").Value = 4670.24
$ws.Range("I40").Value = 4307.2144
$ws.Range("J40").Value = 5132.273
$ws.Range("K40").Value = 4307.2144
$ws.Range("L40").Value = 5132.273
$ws.Range("M40").Value = -4171.2144
$ws.Range("N40").Value = -5404.273
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 1702.8235
$ws.Range("I46").Value = 911.38464
$ws.Range("J46").Value = 4275
$ws.Range("K46").Value = 911.38464
$ws.Range("L46").Value = 4275
$ws.Range("M46").Value = -723.38464
$ws.Range("N46").Value = -4651
$ws.Range("H55").Value = 58824264
$ws.Range("I55").Value = 90909770
$ws.Range("J55").Value = 840.8333
$ws.Range("K55").Value = 90909770
$ws.Range("L55").Value = 840.8333
$ws.Range("M55").Value = -90909597
$ws.Range("N55").Value = -1186.8333
$ws.Range("H74").Value = 83497.5
$ws.Range("I74").Value = 85000
$ws.Range("K74").Value = 85000
$ws.Range("M74").Value = -84002
$ws.Range("H77").Value = 83497.5
$ws.Range("I77").Value = 85000
$ws.Range("K77").Value = 255000
$ws.Range("M77").Value = -250008
$ws.Range("H93").Value = 43479484
$ws.Range("I93").Value = 71429390
$ws.Range("K93").Value = 71429390
$ws.Range("M93").Value = -71428142
$ws.Range("H100").Value = 3403.5386
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 4500
$ws.Range("N100").Value = -5582
$ws.Range("H126").Value = 10838.571
$ws.Range("I126").Value = 4357.4546
$ws.Range("K126").Value = 13072.3638
$ws.Range("M126").Value = -10602.3638
$ws.Range("H132").Value = 736814.2
$ws.Range("I132").Value = 145543.28
$ws.Range("J132").Value = 1254176.2
$ws.Range("K132").Value = 436629.84
$ws.Range("L132").Value = 3762528.6
$ws.Range("M132").Value = -434099.84
$ws.Range("N132").Value = -3767588.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 90019
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H80").Value = 42432.5
$ws.Range("J80").Value = 54875
$ws.Range("L80").Value = 54875
$ws.Range("N80").Value = -56871
$ws.Range("H83").Value = 42432.5
$ws.Range("J83").Value = 54875
$ws.Range("L83").Value = 164625
$ws.Range("N83").Value = -174609
$ws.Range("H122").Value = 2561.7036
$ws.Range("I122").Value = 2665.5881
$ws.Range("K122").Value = 7996.7643
$ws.Range("M122").Value = -5546.7643
$ws.Range("H132").Value = 296642.7
$ws.Range("I132").Value = 2493.625
$ws.Range("K132").Value = 7480.875
$ws.Range("M132").Value = -5050.875
